# ID: OMS-RTM-01 --> Update RTM approved
# Updates the Requirements Traceability Matrix content (rows 6-14) to the
# newly-approved wording, and resizes the affected rows so the (now
# longer/shorter) wrapped text still fits the same way it does in the
# reference workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Admin login / admin profile requirements -----------------------
$ws.Range("A6").Value = "OMS_SIQ_admin_01 & OMS_SIQ_admin_02 &  OMS_CUST_client_07"
$ws.Range("B6").Value = "OMS_SRS_AdmP_01& OMS_SRS_AdmP_02  & OMS_SRS_AdmP_03 & OMS_SRS_AdmP_04 "
$ws.Rows.Item(6).RowHeight = 70.5

# --- Row 7: Item details requirement (B7 collapsed to single Core-B_03) ----
$ws.Range("B7").Value = " OMS_SRS_Core-B_03 "

# --- Row 8: History requirements (added hist_10 / Hist_02) -----------------
$ws.Range("A8").Value = " OMS_SIQ_hist_10 & OMS_SIQ_hist_11   &  OMS_CUST_client_03"
$ws.Range("B8").Value = "OMS_SRS_ Hist_01 & OMS_SRS_ Hist_02"
$ws.Rows.Item(8).RowHeight = 52.5

# --- Rows 9-11: unchanged text, just re-measured row heights ---------------
$ws.Rows.Item(9).RowHeight = 20.25
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 21

# --- Row 12: client_07 duplicate corrected to client_08 ---------------------
$ws.Range("A12").Value = "OMS_SIQ_client_08"
$ws.Rows.Item(12).RowHeight = 19.5

# --- Row 13: client_08 -> client_09, SRS mapping updated to Core-B_07.x -----
$ws.Range("A13").Value = "OMS_SIQ_client_09"
$ws.Range("B13").Value = "OMS_SRS_Core-B_07.01 & OMS_SRS_Core-B_07.02 & OMS_SRS_Core-B_07.03"
$ws.Rows.Item(13).RowHeight = 51

# --- Row 14: new combined customer requirement row + Core-B_08 -------------
$ws.Range("A14").Value = "OMS_CUST_client_01 & OMS_CUST_client_02 & OMS_CUST_client_03 & OMS_CUST_client_04 & OMS_CUST_client_05  &  OMS_CUST_client_06 & OMS_CUST_client_07 &OMS_CUST_client_08"
$ws.Range("B14").Value = "OMS_SRS_Core-B_08"
$ws.Rows.Item(14).RowHeight = 111

# --- Rows 219/220: align with surrounding 15.75pt rows ----------------------
$ws.Rows.Item(219).RowHeight = 15.75
$ws.Rows.Item(220).RowHeight = 15.75
